$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.258.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = "'1.677.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.74%  '
$ws.Range("D5").Value = "'212.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.31%  '
$ws.Range("D6").Value = "'0.5262"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.95%  '
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("D9").Value = "'0.06297"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.52%  '
$ws.Range("D10").Value = "'21.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.05%  '
$ws.Range("D11").Value = "'0.07556"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").Value = "'1.680.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("D13").Value = "'4.472"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("D14").Value = "'0.5636"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.46%  '
$ws.Range("D15").Value = "'67.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.000008046"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.21%  '
$ws.Range("D17").Value = "'26.024.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.73%  '
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").Value = "'4.830"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.50%  '
$ws.Range("D20").Value = "'188.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").Value = "'10.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.03%  '
$ws.Range("D22").Value = "'6.200"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.80%  '
$ws.Range("D24").Value = "'150.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("D25").Value = "'0.1258"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.40%  '
$ws.Range("D26").Value = "'7.599"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.90%  '
$ws.Range("D27").Value = "'16.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("D28").Value = "'0.06223"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("D29").Value = "'1.361"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.39%  '
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("D31").Value = "'3.510"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.94%  '
$ws.Range("D32").Value = "'3.449"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.41%  '
$ws.Range("D33").Value = "'1.635"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.29%  '
$ws.Range("D34").Value = "'1.004"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.62%  '
$ws.Range("D35").Value = "'0.6076"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.67%  '
$ws.Range("D36").Value = "'2.406"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").Value = "'2.737"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("D38").Value = "'0.01622"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.26%  '
$ws.Range("D39").Value = "'6.105"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("D40").Value = "'1.104.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("D41").Value = "'0.8717"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("D43").Value = "'100.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.23%  '
$ws.Range("D44").Value = "'1.825.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("D45").Value = "'0.00000000108"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").Value = "'56.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.50%  '
$ws.Range("D47").Value = "'1.005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D48").Value = "'8.054"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.59%  '
$ws.Range("D49").Value = "'0.05232"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").Value = "'0.4255"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.16%  '
$ws.Range("D51").Value = "'5.988"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.81%  '
